# Insert two new weekly price rows for Brócoli (Vega Monumental Concepción)
# dated 2022-03-18 (serial 44637), pushing all existing rows from 210
# onward down by two (210-241 -> 212-243).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("210:211").Insert()

# Row 210: Primera
$ws.Range("A210").Value = 11
$ws.Range("B210").Value = "Vega Monumental Concepción"
$ws.Range("C210").Value = "Bíobío"
$ws.Range("D210").Value = 44637
$ws.Range("E210").Value = 8
$ws.Range("F210").Value = 100112023
$ws.Range("G210").Value = "Brócoli"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 1000
$ws.Range("K210").Value = 1000
$ws.Range("L210").Value = 1000
$ws.Range("M210").Value = 1000
$ws.Range("N210").Value = "$/unidad"
$ws.Range("O210").Value = "Región Metropolitana"
$ws.Range("P210").Value = 1000
$ws.Range("Q210").Value = 1
$ws.Range("R210").Value = "Hortaliza"

# Row 211: Segunda
$ws.Range("A211").Value = 11
$ws.Range("B211").Value = "Vega Monumental Concepción"
$ws.Range("C211").Value = "Bíobío"
$ws.Range("D211").Value = 44637
$ws.Range("E211").Value = 8
$ws.Range("F211").Value = 100112023
$ws.Range("G211").Value = "Brócoli"
$ws.Range("H211").Value = "Sin especificar"
$ws.Range("I211").Value = "Segunda"
$ws.Range("J211").Value = 1200
$ws.Range("K211").Value = 800
$ws.Range("L211").Value = 800
$ws.Range("M211").Value = 800
$ws.Range("N211").Value = "$/unidad"
$ws.Range("O211").Value = "Región Metropolitana"
$ws.Range("P211").Value = 800
$ws.Range("Q211").Value = 1
$ws.Range("R211").Value = "Hortaliza"
